$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.392.57"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "2.409.56"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("D9").Value = "2.445.83"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.321"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.23%  "
$ws.Range("D14").Value = "2.838.90"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").Value = "57.452.11"
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "2.452.95"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "2.541.67"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("E28").Value = "  -5.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.153"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("E39").Value = "  +3.95%  "
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "134.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "256.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("E51").Value = "  +0.58%  "
